# New submission synced into the "JSS 3D" sheet's results table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3D")

# Append the new response as row 13 (the row right after the last existing
# data row, 12), matching the columns: Timestamp, Full Name, Admission No, AI Score.
$ws.Cells.Item(13, 1).Value = "2026-02-12 19:51:17"
$ws.Cells.Item(13, 2).Value = "Hauwa Hussaini maina"

# "Admission No" is a text column even though this value looks numeric
# (e.g. row 12 stores "19" as text, not a number) -- use a leading
# apostrophe so Excel keeps it as text instead of auto-converting to a
# number, then reset the style so the quote-prefix formatting marker
# left behind doesn't linger on the cell.
$ws.Cells.Item(13, 3).Value = "'40"
$ws.Cells.Item(13, 3).Style = "Normal"

$ws.Cells.Item(13, 4).Value = 7
